$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.6
$ws.Range("G2").Value = 1.62
$ws.Range("H2").Value = 6.2
$ws.Range("I2").Value = 6.6
$ws.Range("P2").Value = 2.14
$ws.Range("Q2").Value = 1.84
$ws.Range("R2").Value = 1.44
$ws.Range("S2").Value = 3.2
$ws.Range("U2").Value = 2.06
$ws.Range("V2").Value = 1.18
$ws.Range("W2").Value = 2.6
$ws.Range("Y2").Value = 23
$ws.Range("AA2").Value = 180
$ws.Range("AD2").Value = 24
$ws.Range("AH2").Value = 22
$ws.Range("AJ2").Value = 15
$ws.Range("AL2").Value = 34
$ws.Range("AO2").Value = 95
$ws.Range("AD3").Value = 14.5
$ws.Range("G4").Value = 3.15
$ws.Range("I4").Value = 3
$ws.Range("K4").Value = 4.2
$ws.Range("M4").Value = 1.08
$ws.Range("V4").Value = 1.52
$ws.Range("W4").Value = 1.47
$ws.Range("O5").Value = 1.24
$ws.Range("P5").Value = 2.22
$ws.Range("Q5").Value = 1.62
$ws.Range("U5").Value = 1.93
$ws.Range("Y5").Value = 9.800000000000001
$ws.Range("AD5").Value = 10
$ws.Range("AE5").Value = 1000
$ws.Range("AN5").Value = 130
$ws.Range("AO5").Value = 7
$ws.Range("F6").Value = 1.4
$ws.Range("G6").Value = 1.45
$ws.Range("K6").Value = 5.5
$ws.Range("L6").Value = 1.33
$ws.Range("N6").Value = 4.3
$ws.Range("P6").Value = 2.18
$ws.Range("Q6").Value = 1.74
$ws.Range("U6").Value = 1.86
$ws.Range("Y6").Value = 30
$ws.Range("AB6").Value = 10
$ws.Range("AE6").Value = 160
$ws.Range("AG6").Value = 12
$ws.Range("AK6").Value = 15
$ws.Range("AL6").Value = 38
$ws.Range("AN6").Value = 6.8
$ws.Range("H7").Value = 4.8
$ws.Range("L7").Value = 1.41
$ws.Range("N7").Value = 3.75
$ws.Range("P7").Value = 1.9
$ws.Range("Q7").Value = 1.97
$ws.Range("X7").Value = 1000
$ws.Range("AC7").Value = 1000
$ws.Range("AK7").Value = 23
$ws.Range("H8").Value = 11
$ws.Range("I8").Value = 11.5
$ws.Range("J8").Value = 4.9
$ws.Range("K8").Value = 5
$ws.Range("L8").Value = 1.39
$ws.Range("Q8").Value = 1.93
$ws.Range("R8").Value = 1.41
$ws.Range("T8").Value = 2.24
$ws.Range("V8").Value = 1.09
$ws.Range("X8").Value = 18.5
$ws.Range("Y8").Value = 32
$ws.Range("AA8").Value = 450
$ws.Range("AE8").Value = 200
$ws.Range("AG8").Value = 10.5
$ws.Range("AM8").Value = 220
$ws.Range("AN8").Value = 7
$ws.Range("J9").Value = 3.25
$ws.Range("K9").Value = 3.55
$ws.Range("L9").Value = 1.45
$ws.Range("P9").Value = 1.82
$ws.Range("Q9").Value = 2.12
$ws.Range("R9").Value = 1.29
$ws.Range("S9").Value = 3.85
$ws.Range("V9").Value = 1.38
$ws.Range("AF9").Value = 17
$ws.Range("F10").Value = 3.85
$ws.Range("I10").Value = 2.24
$ws.Range("L10").Value = 1.46
$ws.Range("P10").Value = 1.79
$ws.Range("Q10").Value = 2.12
$ws.Range("U10").Value = 2.02
$ws.Range("V10").Value = 1.8
$ws.Range("AA10").Value = 34
$ws.Range("AC10").Value = 9
$ws.Range("AD10").Value = 12.5
$ws.Range("AE10").Value = 30
$ws.Range("AF10").Value = 34
$ws.Range("AH10").Value = 21
$ws.Range("H11").Value = 1.61
$ws.Range("I11").Value = 1.62
$ws.Range("K11").Value = 4.1
$ws.Range("P11").Value = 1.81
$ws.Range("R11").Value = 1.3
$ws.Range("U11").Value = 1.78
$ws.Range("V11").Value = 2.6
$ws.Range("AB11").Value = 20
$ws.Range("AL11").Value = 130
$ws.Range("AM11").Value = 190
$ws.Range("AN11").Value = 200
$ws.Range("F12").Value = 2.86
$ws.Range("H12").Value = 2.88
$ws.Range("I12").Value = 2.9
$ws.Range("V12").Value = 1.52
$ws.Range("X12").Value = 9.800000000000001
$ws.Range("Y12").Value = 9
$ws.Range("Z12").Value = 16.5
$ws.Range("AC12").Value = 7.2
$ws.Range("AF12").Value = 16.5
$ws.Range("AG12").Value = 13.5
$ws.Range("AN12").Value = 40
$ws.Range("G13").Value = 4.7
$ws.Range("H13").Value = 2.02
$ws.Range("J13").Value = 3.4
$ws.Range("P13").Value = 1.72
$ws.Range("Q13").Value = 2.26
$ws.Range("R13").Value = 1.26
$ws.Range("T13").Value = 1.99
$ws.Range("U13").Value = 1.93
$ws.Range("W13").Value = 1.27
$ws.Range("X13").Value = 10.5
$ws.Range("Z13").Value = 11.5
$ws.Range("AA13").Value = 24
$ws.Range("AC13").Value = 7.6
$ws.Range("AD13").Value = 10.5
$ws.Range("AE13").Value = 24
$ws.Range("AF13").Value = 30
$ws.Range("AG13").Value = 18
$ws.Range("AI13").Value = 46
$ws.Range("AJ13").Value = 110
$ws.Range("AK13").Value = 65
$ws.Range("AM13").Value = 140
$ws.Range("AN13").Value = 95
$ws.Range("AO13").Value = 19.5
$ws.Range("F14").Value = 3.25
$ws.Range("H14").Value = 2.56
$ws.Range("I14").Value = 2.66
$ws.Range("N14").Value = 2.74
$ws.Range("P14").Value = 1.59
$ws.Range("Q14").Value = 2.56
$ws.Range("S14").Value = 5.3
$ws.Range("T14").Value = 2.08
$ws.Range("W14").Value = 1.42
$ws.Range("X14").Value = 8.800000000000001
$ws.Range("AO14").Value = 38
$ws.Range("F15").Value = 1.91
$ws.Range("G15").Value = 1.95
$ws.Range("H15").Value = 4.1
$ws.Range("J15").Value = 3.8
$ws.Range("P15").Value = 2.14
$ws.Range("W15").Value = 2.04
$ws.Range("Z15").Value = 980
$ws.Range("J16").Value = 3.45
$ws.Range("K16").Value = 3.5
$ws.Range("N16").Value = 3.85
$ws.Range("O16").Value = 1.33
$ws.Range("Y16").Value = 13
$ws.Range("AD16").Value = 15
$ws.Range("AK16").Value = 24
$ws.Range("AM16").Value = 85
$ws.Range("AN16").Value = 18.5
$ws.Range("G17").Value = 1.87
$ws.Range("H17").Value = 5.1
$ws.Range("I17").Value = 5.2
$ws.Range("K17").Value = 3.75
$ws.Range("T17").Value = 1.95
$ws.Range("AD17").Value = 20
$ws.Range("F18").Value = 4
$ws.Range("I18").Value = 2.18
$ws.Range("M18").Value = 1.09
$ws.Range("Q18").Value = 2.12
$ws.Range("T18").Value = 1.87
$ws.Range("U18").Value = 2.06
$ws.Range("V18").Value = 1.85
$ws.Range("X18").Value = 24
$ws.Range("Y18").Value = 8.800000000000001
$ws.Range("AH18").Value = 20
$ws.Range("X19").Value = 16
$ws.Range("AH19").Value = 16
$ws.Range("AM19").Value = 95
$ws.Range("L20").Value = 1.18
$ws.Range("R20").Value = 1.86
$ws.Range("T20").Value = 2.38
$ws.Range("U20").Value = 1.68
$ws.Range("F21").Value = 1.5
$ws.Range("M21").Value = 1.05
